$d = $word.ActiveDocument

# Remove all Custom XML Parts (e.g. SharePoint content-type schemas and
# document-library form metadata) that were injected into the package.
# This mirrors the diff that drops customXml/item1.xml, item2.xml, item3.xml
# and their accompanying itemProps*.xml parts from the OOXML package.
$count = $d.CustomXMLParts.Count
for ($i = $count; $i -ge 1; $i--) {
    $part = $d.CustomXMLParts.Item($i)
    $part.Delete()
}
